# Edit script: add accessibility requirements (section 2.x) to requirements.docx
# and tidy up the 1.2 requirement text (merge the proofread-marked runs into one run).
#
# NOTE: Paragraph.Index is not a reliable "absolute position in the document"
# indicator in this runtime (it appears to reset across stories/tables), so
# all paragraph addressing below uses plain 1-based position within
# $d.Paragraphs, tracked manually.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function New-PkgXml($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + ' ' + $w14Ns + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Find-ParaPosition($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Paragraph "1.2 ..." — collapse the proofed ("be clarified" / " in section 2")
#    runs back into a single plain run.
# ---------------------------------------------------------------------------

$pos12 = Find-ParaPosition("1.2 The User shall be able to delete books")

$p12Body = '<w:p w14:paraId="17F3B1EB" w14:textId="4278DA13" w:rsidR="006A01A0" w:rsidRDefault="006A01A0" w:rsidP="006A01A0"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr></w:pPr><w:r w:rsidRPr="006A01A0"><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr><w:t>1.2 The User shall be able to delete books from the app. This will delete the copy of the book and any notes/bookmarks that the book has, notes/bookmarks will be clarified in section 2</w:t></w:r></w:p>'

$d.Paragraphs.Item($pos12).Range.InsertXML((New-PkgXml $p12Body))

# ---------------------------------------------------------------------------
# 2) After paragraph "1.5 ..." insert: a blank line, then the new "2." block
#    (2, 2.1, 2.2 in one paragraph separated by line breaks), then "2.3", then
#    "2.4" (which itself has a proofread-marked "similar to").
# ---------------------------------------------------------------------------

$pos15 = Find-ParaPosition("1.5 The User shall be able to remove categories")

$emptyBody = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/></w:rPr></w:pPr></w:p>'

$sec2Body = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>2. The User shall be able to make viewing adjustments for accessibility purposes</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:br/><w:t>2.1 The User shall be able to swipe across screen right to left to move to the next page and left to right to return to the previous page</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:br/><w:t>2.2 The User shall be able to mark a page and be able to jump to that marked page immediately</w:t></w:r></w:p>'

$sec23Body = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>2.3 The User shall be able to switch between Night mode and Day mode</w:t></w:r></w:p>'

$sec24Body = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t xml:space="preserve">2.4 The User shall be able to search for a word and if the word is found they shall be able to jump to it in a fashion </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t>similar to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:eastAsia="en-US"/></w:rPr><w:t xml:space="preserve"> a jump to the marked page above</w:t></w:r></w:p>'

# Insert the four new paragraphs, in order, right after "1.5 ...".
# Each InsertParagraphAfter() creates a fresh empty paragraph immediately
# after the anchor; we then overwrite its contents via InsertXML. We track
# the insertion position manually (plain integer), since Paragraph.Index is
# unreliable in this runtime.

$pos = $pos15

$anchor = $d.Paragraphs.Item($pos).Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pos = $pos + 1
$d.Paragraphs.Item($pos).Range.InsertXML((New-PkgXml $emptyBody))

$anchor = $d.Paragraphs.Item($pos).Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pos = $pos + 1
$d.Paragraphs.Item($pos).Range.InsertXML((New-PkgXml $sec2Body))

$anchor = $d.Paragraphs.Item($pos).Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pos = $pos + 1
$d.Paragraphs.Item($pos).Range.InsertXML((New-PkgXml $sec23Body))

$anchor = $d.Paragraphs.Item($pos).Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$pos = $pos + 1
$d.Paragraphs.Item($pos).Range.InsertXML((New-PkgXml $sec24Body))

Write-Output "done"
